$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number-format strings matching cellXfs already present in xl/styles.xml.
# Re-using the exact same format-code text makes the engine map each cell to
# the pre-existing style index (2/3/4/5) instead of minting new ones.
$fmtDate = "dd\.mm\.yyyy"
$fmtInt  = "###############"
$fmtNum  = "#,###.00"

# Row 218
$ws.Range("A218").NumberFormat = $fmtDate
$ws.Range("A218").Value = (Get-Date -Year 2025 -Month 6 -Day 1 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B218").HorizontalAlignment = -4152
$ws.Range("B218").Value = 'שיק'

$ws.Range("D218").NumberFormat = $fmtInt
$ws.Range("D218").Value = 7302259

$ws.Range("E218").NumberFormat = $fmtNum
$ws.Range("E218").Value = 5300

$ws.Range("F218").NumberFormat = $fmtNum
$ws.Range("F218").Value = ""

$ws.Range("G218").NumberFormat = $fmtNum
$ws.Range("G218").Value = 28835.43

$ws.Range("H218").NumberFormat = $fmtDate
$ws.Range("H218").Value = (Get-Date -Year 2025 -Month 6 -Day 1 -Hour 0 -Minute 0 -Second 0)

# Row 219
$ws.Range("A219").NumberFormat = $fmtDate
$ws.Range("A219").Value = (Get-Date -Year 2025 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B219").HorizontalAlignment = -4152
$ws.Range("B219").Value = 'מסטרקרד'

$ws.Range("D219").NumberFormat = $fmtInt
$ws.Range("D219").Value = 961

$ws.Range("E219").NumberFormat = $fmtNum
$ws.Range("E219").Value = 2929.31

$ws.Range("F219").NumberFormat = $fmtNum
$ws.Range("F219").Value = ""

$ws.Range("G219").NumberFormat = $fmtNum
$ws.Range("G219").Value = 31206.12

$ws.Range("H219").NumberFormat = $fmtDate
$ws.Range("H219").Value = (Get-Date -Year 2025 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)

# Row 220
$ws.Range("A220").NumberFormat = $fmtDate
$ws.Range("A220").Value = (Get-Date -Year 2025 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B220").HorizontalAlignment = -4152
$ws.Range("B220").Value = 'החזרת שיק-טכני (תאריך ערך: 01/06)'

$ws.Range("C220").HorizontalAlignment = -4152
$ws.Range("C220").Value = ' סיבת החזרה:  פעולה אסורה עפ"י דין'

$ws.Range("D220").NumberFormat = $fmtInt
$ws.Range("D220").Value = 7302259

$ws.Range("E220").NumberFormat = $fmtNum
$ws.Range("E220").Value = ""

$ws.Range("F220").NumberFormat = $fmtNum
$ws.Range("F220").Value = 5300

$ws.Range("G220").NumberFormat = $fmtNum
$ws.Range("G220").Value = 34135.43

$ws.Range("H220").NumberFormat = $fmtDate
$ws.Range("H220").Value = (Get-Date -Year 2025 -Month 6 -Day 1 -Hour 0 -Minute 0 -Second 0)

$ws.Range("I220").HorizontalAlignment = -4152
$ws.Range("I220").Value = 'סיבת החזרה:'

$ws.Range("J220").HorizontalAlignment = -4152
$ws.Range("J220").Value = 'פעולה אסורה עפ"י דין'

# Row 221
$ws.Range("A221").NumberFormat = $fmtDate
$ws.Range("A221").Value = (Get-Date -Year 2025 -Month 6 -Day 4 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B221").HorizontalAlignment = -4152
$ws.Range("B221").Value = 'העב'' לאחר-נייד'

$ws.Range("C221").HorizontalAlignment = -4152
$ws.Range("C221").Value = 'לטובת: פני דויטש עבור: שכירות יוני 2025'

$ws.Range("D221").NumberFormat = $fmtInt
$ws.Range("D221").Value = 399108067

$ws.Range("E221").NumberFormat = $fmtNum
$ws.Range("E221").Value = 5300

$ws.Range("F221").NumberFormat = $fmtNum
$ws.Range("F221").Value = ""

$ws.Range("G221").NumberFormat = $fmtNum
$ws.Range("G221").Value = 29096.12

$ws.Range("H221").NumberFormat = $fmtDate
$ws.Range("H221").Value = (Get-Date -Year 2025 -Month 6 -Day 4 -Hour 0 -Minute 0 -Second 0)

$ws.Range("I221").HorizontalAlignment = -4152
$ws.Range("I221").Value = 'פני דויטש'

$ws.Range("J221").HorizontalAlignment = -4152
$ws.Range("J221").Value = 'שכירות יוני 2025'

# Row 222
$ws.Range("A222").NumberFormat = $fmtDate
$ws.Range("A222").Value = (Get-Date -Year 2025 -Month 6 -Day 4 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B222").HorizontalAlignment = -4152
$ws.Range("B222").Value = 'ע. מסלול בסיסי'

$ws.Range("D222").NumberFormat = $fmtInt
$ws.Range("D222").Value = 817

$ws.Range("E222").NumberFormat = $fmtNum
$ws.Range("E222").Value = 10

$ws.Range("F222").NumberFormat = $fmtNum
$ws.Range("F222").Value = ""

$ws.Range("G222").NumberFormat = $fmtNum
$ws.Range("G222").Value = 34396.120000000003

$ws.Range("H222").NumberFormat = $fmtDate
$ws.Range("H222").Value = (Get-Date -Year 2025 -Month 6 -Day 4 -Hour 0 -Minute 0 -Second 0)

# Row 223
$ws.Range("A223").NumberFormat = $fmtDate
$ws.Range("A223").Value = (Get-Date -Year 2025 -Month 6 -Day 4 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B223").HorizontalAlignment = -4152
$ws.Range("B223").Value = 'הפק.שיק בסלולר'

$ws.Range("D223").NumberFormat = $fmtInt
$ws.Range("D223").Value = 80000051

$ws.Range("E223").NumberFormat = $fmtNum
$ws.Range("E223").Value = ""

$ws.Range("F223").NumberFormat = $fmtNum
$ws.Range("F223").Value = 3200

$ws.Range("G223").NumberFormat = $fmtNum
$ws.Range("G223").Value = 34406.120000000003

$ws.Range("H223").NumberFormat = $fmtDate
$ws.Range("H223").Value = (Get-Date -Year 2025 -Month 6 -Day 4 -Hour 0 -Minute 0 -Second 0)

# Row 224
$ws.Range("A224").NumberFormat = $fmtDate
$ws.Range("A224").Value = (Get-Date -Year 2025 -Month 6 -Day 9 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B224").HorizontalAlignment = -4152
$ws.Range("B224").Value = 'העב'' לאחר-נייד'

$ws.Range("C224").HorizontalAlignment = -4152
$ws.Range("C224").Value = 'לטובת: SACHIN CHANDRAPPA עבור: salary  May 2025'

$ws.Range("D224").NumberFormat = $fmtInt
$ws.Range("D224").Value = 400034110

$ws.Range("E224").NumberFormat = $fmtNum
$ws.Range("E224").Value = 3493

$ws.Range("F224").NumberFormat = $fmtNum
$ws.Range("F224").Value = ""

$ws.Range("G224").NumberFormat = $fmtNum
$ws.Range("G224").Value = 25603.119999999999

$ws.Range("H224").NumberFormat = $fmtDate
$ws.Range("H224").Value = (Get-Date -Year 2025 -Month 6 -Day 9 -Hour 0 -Minute 0 -Second 0)

$ws.Range("I224").HorizontalAlignment = -4152
$ws.Range("I224").Value = 'SACHIN CHANDRAPPA'

$ws.Range("J224").HorizontalAlignment = -4152
$ws.Range("J224").Value = 'salary  May 2025'

# Row 225
$ws.Range("A225").NumberFormat = $fmtDate
$ws.Range("A225").Value = (Get-Date -Year 2025 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B225").HorizontalAlignment = -4152
$ws.Range("B225").Value = 'שירותי בריאות'

$ws.Range("D225").NumberFormat = $fmtInt
$ws.Range("D225").Value = 30744700

$ws.Range("E225").NumberFormat = $fmtNum
$ws.Range("E225").Value = 72.11

$ws.Range("F225").NumberFormat = $fmtNum
$ws.Range("F225").Value = ""

$ws.Range("G225").NumberFormat = $fmtNum
$ws.Range("G225").Value = 25304.400000000001

$ws.Range("H225").NumberFormat = $fmtDate
$ws.Range("H225").Value = (Get-Date -Year 2025 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Row 226
$ws.Range("A226").NumberFormat = $fmtDate
$ws.Range("A226").Value = (Get-Date -Year 2025 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B226").HorizontalAlignment = -4152
$ws.Range("B226").Value = 'חברת החשמל ליש'

$ws.Range("D226").NumberFormat = $fmtInt
$ws.Range("D226").Value = 706994

$ws.Range("E226").NumberFormat = $fmtNum
$ws.Range("E226").Value = 226.61

$ws.Range("F226").NumberFormat = $fmtNum
$ws.Range("F226").Value = ""

$ws.Range("G226").NumberFormat = $fmtNum
$ws.Range("G226").Value = 25376.51

$ws.Range("H226").NumberFormat = $fmtDate
$ws.Range("H226").Value = (Get-Date -Year 2025 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Row 227
$ws.Range("A227").NumberFormat = $fmtDate
$ws.Range("A227").Value = (Get-Date -Year 2025 -Month 6 -Day 15 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B227").HorizontalAlignment = -4152
$ws.Range("B227").Value = 'הו"ק פקדון'

$ws.Range("D227").NumberFormat = $fmtInt
$ws.Range("D227").Value = 2590001

$ws.Range("E227").NumberFormat = $fmtNum
$ws.Range("E227").Value = 1100

$ws.Range("F227").NumberFormat = $fmtNum
$ws.Range("F227").Value = ""

$ws.Range("G227").NumberFormat = $fmtNum
$ws.Range("G227").Value = 24204.400000000001

$ws.Range("H227").NumberFormat = $fmtDate
$ws.Range("H227").Value = (Get-Date -Year 2025 -Month 6 -Day 15 -Hour 0 -Minute 0 -Second 0)

# Row 228
$ws.Range("A228").NumberFormat = $fmtDate
$ws.Range("A228").Value = (Get-Date -Year 2025 -Month 6 -Day 16 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B228").HorizontalAlignment = -4152
$ws.Range("B228").Value = 'סקיולייף-נטלי'

$ws.Range("D228").NumberFormat = $fmtInt
$ws.Range("D228").Value = 32349623

$ws.Range("E228").NumberFormat = $fmtNum
$ws.Range("E228").Value = 25

$ws.Range("F228").NumberFormat = $fmtNum
$ws.Range("F228").Value = ""

$ws.Range("G228").NumberFormat = $fmtNum
$ws.Range("G228").Value = 24179.4

$ws.Range("H228").NumberFormat = $fmtDate
$ws.Range("H228").Value = (Get-Date -Year 2025 -Month 6 -Day 16 -Hour 0 -Minute 0 -Second 0)

# Row 229
$ws.Range("A229").NumberFormat = $fmtDate
$ws.Range("A229").Value = (Get-Date -Year 2025 -Month 6 -Day 20 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B229").HorizontalAlignment = -4152
$ws.Range("B229").Value = 'בזק-הוראות קבע'

$ws.Range("D229").NumberFormat = $fmtInt
$ws.Range("D229").Value = 703991

$ws.Range("E229").NumberFormat = $fmtNum
$ws.Range("E229").Value = 49.16

$ws.Range("F229").NumberFormat = $fmtNum
$ws.Range("F229").Value = ""

$ws.Range("G229").NumberFormat = $fmtNum
$ws.Range("G229").Value = 24130.240000000002

$ws.Range("H229").NumberFormat = $fmtDate
$ws.Range("H229").Value = (Get-Date -Year 2025 -Month 6 -Day 20 -Hour 0 -Minute 0 -Second 0)

# Row 230
$ws.Range("A230").NumberFormat = $fmtDate
$ws.Range("A230").Value = (Get-Date -Year 2025 -Month 6 -Day 25 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B230").HorizontalAlignment = -4152
$ws.Range("B230").Value = 'ביטוח.ל-אזו"ש'

$ws.Range("C230").HorizontalAlignment = -4152
$ws.Range("C230").Value = '  עבור: מימון אסתר          מזהה 070120555'

$ws.Range("D230").NumberFormat = $fmtInt
$ws.Range("D230").Value = 13101084

$ws.Range("E230").NumberFormat = $fmtNum
$ws.Range("E230").Value = ""

$ws.Range("F230").NumberFormat = $fmtNum
$ws.Range("F230").Value = 1422

$ws.Range("G230").NumberFormat = $fmtNum
$ws.Range("G230").Value = 27937.24

$ws.Range("H230").NumberFormat = $fmtDate
$ws.Range("H230").Value = (Get-Date -Year 2025 -Month 6 -Day 25 -Hour 0 -Minute 0 -Second 0)

$ws.Range("J230").HorizontalAlignment = -4152
$ws.Range("J230").Value = 'מימון אסתר          מזהה 070120555'

# Row 231
$ws.Range("A231").NumberFormat = $fmtDate
$ws.Range("A231").Value = (Get-Date -Year 2025 -Month 6 -Day 25 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B231").HorizontalAlignment = -4152
$ws.Range("B231").Value = 'ביטוח.ל-אזו"ש'

$ws.Range("C231").HorizontalAlignment = -4152
$ws.Range("C231").Value = '  עבור: מימון אסתר          מזהה 070120555'

$ws.Range("D231").NumberFormat = $fmtInt
$ws.Range("D231").Value = 13101084

$ws.Range("E231").NumberFormat = $fmtNum
$ws.Range("E231").Value = ""

$ws.Range("F231").NumberFormat = $fmtNum
$ws.Range("F231").Value = 2385

$ws.Range("G231").NumberFormat = $fmtNum
$ws.Range("G231").Value = 26515.24

$ws.Range("H231").NumberFormat = $fmtDate
$ws.Range("H231").Value = (Get-Date -Year 2025 -Month 6 -Day 25 -Hour 0 -Minute 0 -Second 0)

$ws.Range("J231").HorizontalAlignment = -4152
$ws.Range("J231").Value = 'מימון אסתר          מזהה 070120555'

# Row 232
$ws.Range("A232").NumberFormat = $fmtDate
$ws.Range("A232").Value = (Get-Date -Year 2025 -Month 6 -Day 30 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B232").HorizontalAlignment = -4152
$ws.Range("B232").Value = 'האוצר-תגמולים'

$ws.Range("C232").HorizontalAlignment = -4152
$ws.Range("C232").Value = '  עבור: מימון אסתר          מזהה 070120555'

$ws.Range("D232").NumberFormat = $fmtInt
$ws.Range("D232").Value = 4061016

$ws.Range("E232").NumberFormat = $fmtNum
$ws.Range("E232").Value = ""

$ws.Range("F232").NumberFormat = $fmtNum
$ws.Range("F232").Value = 12490.45

$ws.Range("G232").NumberFormat = $fmtNum
$ws.Range("G232").Value = 40427.69

$ws.Range("H232").NumberFormat = $fmtDate
$ws.Range("H232").Value = (Get-Date -Year 2025 -Month 6 -Day 30 -Hour 0 -Minute 0 -Second 0)

$ws.Range("J232").HorizontalAlignment = -4152
$ws.Range("J232").Value = 'מימון אסתר          מזהה 070120555'

# Row 233
$ws.Range("A233").NumberFormat = $fmtDate
$ws.Range("A233").Value = (Get-Date -Year 2025 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B233").HorizontalAlignment = -4152
$ws.Range("B233").Value = 'מסטרקרד'

$ws.Range("D233").NumberFormat = $fmtInt
$ws.Range("D233").Value = 961

$ws.Range("E233").NumberFormat = $fmtNum
$ws.Range("E233").Value = 1922.17

$ws.Range("F233").NumberFormat = $fmtNum
$ws.Range("F233").Value = ""

$ws.Range("G233").NumberFormat = $fmtNum
$ws.Range("G233").Value = 33205.519999999997

$ws.Range("H233").NumberFormat = $fmtDate
$ws.Range("H233").Value = (Get-Date -Year 2025 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0)

# Row 234
$ws.Range("A234").NumberFormat = $fmtDate
$ws.Range("A234").Value = (Get-Date -Year 2025 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B234").HorizontalAlignment = -4152
$ws.Range("B234").Value = 'העב'' לאחר-נייד'

$ws.Range("C234").HorizontalAlignment = -4152
$ws.Range("C234").Value = 'לטובת: פני דויטש עבור: שכירות יולי'

$ws.Range("D234").NumberFormat = $fmtInt
$ws.Range("D234").Value = 403621594

$ws.Range("E234").NumberFormat = $fmtNum
$ws.Range("E234").Value = 5300

$ws.Range("F234").NumberFormat = $fmtNum
$ws.Range("F234").Value = ""

$ws.Range("G234").NumberFormat = $fmtNum
$ws.Range("G234").Value = 35127.69

$ws.Range("H234").NumberFormat = $fmtDate
$ws.Range("H234").Value = (Get-Date -Year 2025 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0)

$ws.Range("I234").HorizontalAlignment = -4152
$ws.Range("I234").Value = 'פני דויטש'

$ws.Range("J234").HorizontalAlignment = -4152
$ws.Range("J234").Value = 'שכירות יולי'

# Row 235
$ws.Range("A235").NumberFormat = $fmtDate
$ws.Range("A235").Value = (Get-Date -Year 2025 -Month 7 -Day 3 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B235").HorizontalAlignment = -4152
$ws.Range("B235").Value = 'ע.מפעולות-ישיר'

$ws.Range("D235").NumberFormat = $fmtInt
$ws.Range("D235").Value = 2

$ws.Range("E235").NumberFormat = $fmtNum
$ws.Range("E235").Value = 3.5

$ws.Range("F235").NumberFormat = $fmtNum
$ws.Range("F235").Value = ""

$ws.Range("G235").NumberFormat = $fmtNum
$ws.Range("G235").Value = 33192.019999999997

$ws.Range("H235").NumberFormat = $fmtDate
$ws.Range("H235").Value = (Get-Date -Year 2025 -Month 7 -Day 3 -Hour 0 -Minute 0 -Second 0)

# Row 236
$ws.Range("A236").NumberFormat = $fmtDate
$ws.Range("A236").Value = (Get-Date -Year 2025 -Month 7 -Day 3 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B236").HorizontalAlignment = -4152
$ws.Range("B236").Value = 'ע. מסלול בסיסי'

$ws.Range("D236").NumberFormat = $fmtInt
$ws.Range("D236").Value = 817

$ws.Range("E236").NumberFormat = $fmtNum
$ws.Range("E236").Value = 10

$ws.Range("F236").NumberFormat = $fmtNum
$ws.Range("F236").Value = ""

$ws.Range("G236").NumberFormat = $fmtNum
$ws.Range("G236").Value = 33195.519999999997

$ws.Range("H236").NumberFormat = $fmtDate
$ws.Range("H236").Value = (Get-Date -Year 2025 -Month 7 -Day 3 -Hour 0 -Minute 0 -Second 0)

# Row 237
$ws.Range("A237").NumberFormat = $fmtDate
$ws.Range("A237").Value = (Get-Date -Year 2025 -Month 7 -Day 6 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B237").HorizontalAlignment = -4152
$ws.Range("B237").Value = 'הפק.שיק בסלולר'

$ws.Range("D237").NumberFormat = $fmtInt
$ws.Range("D237").Value = 80000052

$ws.Range("E237").NumberFormat = $fmtNum
$ws.Range("E237").Value = ""

$ws.Range("F237").NumberFormat = $fmtNum
$ws.Range("F237").Value = 3200

$ws.Range("G237").NumberFormat = $fmtNum
$ws.Range("G237").Value = 36392.019999999997

$ws.Range("H237").NumberFormat = $fmtDate
$ws.Range("H237").Value = (Get-Date -Year 2025 -Month 7 -Day 6 -Hour 0 -Minute 0 -Second 0)

# Row 238
$ws.Range("A238").NumberFormat = $fmtDate
$ws.Range("A238").Value = (Get-Date -Year 2025 -Month 7 -Day 9 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B238").HorizontalAlignment = -4152
$ws.Range("B238").Value = 'העב'' לאחר-נייד'

$ws.Range("C238").HorizontalAlignment = -4152
$ws.Range("C238").Value = 'לטובת: SACHIN CHANDRAPPA עבור: salary 062025'

$ws.Range("D238").NumberFormat = $fmtInt
$ws.Range("D238").Value = 405155007

$ws.Range("E238").NumberFormat = $fmtNum
$ws.Range("E238").Value = 3314

$ws.Range("F238").NumberFormat = $fmtNum
$ws.Range("F238").Value = ""

$ws.Range("G238").NumberFormat = $fmtNum
$ws.Range("G238").Value = 33078.019999999997

$ws.Range("H238").NumberFormat = $fmtDate
$ws.Range("H238").Value = (Get-Date -Year 2025 -Month 7 -Day 9 -Hour 0 -Minute 0 -Second 0)

$ws.Range("I238").HorizontalAlignment = -4152
$ws.Range("I238").Value = 'SACHIN CHANDRAPPA'

$ws.Range("J238").HorizontalAlignment = -4152
$ws.Range("J238").Value = 'salary 062025'

# Row 239
$ws.Range("A239").NumberFormat = $fmtDate
$ws.Range("A239").Value = (Get-Date -Year 2025 -Month 7 -Day 10 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B239").HorizontalAlignment = -4152
$ws.Range("B239").Value = 'שירותי בריאות'

$ws.Range("D239").NumberFormat = $fmtInt
$ws.Range("D239").Value = 30744700

$ws.Range("E239").NumberFormat = $fmtNum
$ws.Range("E239").Value = 220.8

$ws.Range("F239").NumberFormat = $fmtNum
$ws.Range("F239").Value = ""

$ws.Range("G239").NumberFormat = $fmtNum
$ws.Range("G239").Value = 32857.22

$ws.Range("H239").NumberFormat = $fmtDate
$ws.Range("H239").Value = (Get-Date -Year 2025 -Month 7 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Row 240
$ws.Range("A240").NumberFormat = $fmtDate
$ws.Range("A240").Value = (Get-Date -Year 2025 -Month 7 -Day 15 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B240").HorizontalAlignment = -4152
$ws.Range("B240").Value = 'הו"ק פקדון'

$ws.Range("D240").NumberFormat = $fmtInt
$ws.Range("D240").Value = 2590001

$ws.Range("E240").NumberFormat = $fmtNum
$ws.Range("E240").Value = 1100

$ws.Range("F240").NumberFormat = $fmtNum
$ws.Range("F240").Value = ""

$ws.Range("G240").NumberFormat = $fmtNum
$ws.Range("G240").Value = 31757.22

$ws.Range("H240").NumberFormat = $fmtDate
$ws.Range("H240").Value = (Get-Date -Year 2025 -Month 7 -Day 15 -Hour 0 -Minute 0 -Second 0)

# Row 241
$ws.Range("A241").NumberFormat = $fmtDate
$ws.Range("A241").Value = (Get-Date -Year 2025 -Month 7 -Day 16 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B241").HorizontalAlignment = -4152
$ws.Range("B241").Value = 'סקיולייף-נטלי'

$ws.Range("D241").NumberFormat = $fmtInt
$ws.Range("D241").Value = 32349623

$ws.Range("E241").NumberFormat = $fmtNum
$ws.Range("E241").Value = 25

$ws.Range("F241").NumberFormat = $fmtNum
$ws.Range("F241").Value = ""

$ws.Range("G241").NumberFormat = $fmtNum
$ws.Range("G241").Value = 31732.22

$ws.Range("H241").NumberFormat = $fmtDate
$ws.Range("H241").Value = (Get-Date -Year 2025 -Month 7 -Day 16 -Hour 0 -Minute 0 -Second 0)

# Row 242
$ws.Range("A242").NumberFormat = $fmtDate
$ws.Range("A242").Value = (Get-Date -Year 2025 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B242").HorizontalAlignment = -4152
$ws.Range("B242").Value = 'ביטוח.ל-אזו"ש'

$ws.Range("C242").HorizontalAlignment = -4152
$ws.Range("C242").Value = '  עבור: מימון אסתר          מזהה 070120555'

$ws.Range("D242").NumberFormat = $fmtInt
$ws.Range("D242").Value = 13101084

$ws.Range("E242").NumberFormat = $fmtNum
$ws.Range("E242").Value = ""

$ws.Range("F242").NumberFormat = $fmtNum
$ws.Range("F242").Value = 1422

$ws.Range("G242").NumberFormat = $fmtNum
$ws.Range("G242").Value = 35539.22

$ws.Range("H242").NumberFormat = $fmtDate
$ws.Range("H242").Value = (Get-Date -Year 2025 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)

$ws.Range("J242").HorizontalAlignment = -4152
$ws.Range("J242").Value = 'מימון אסתר          מזהה 070120555'

# Row 243
$ws.Range("A243").NumberFormat = $fmtDate
$ws.Range("A243").Value = (Get-Date -Year 2025 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B243").HorizontalAlignment = -4152
$ws.Range("B243").Value = 'ביטוח.ל-אזו"ש'

$ws.Range("C243").HorizontalAlignment = -4152
$ws.Range("C243").Value = '  עבור: מימון אסתר          מזהה 070120555'

$ws.Range("D243").NumberFormat = $fmtInt
$ws.Range("D243").Value = 13101084

$ws.Range("E243").NumberFormat = $fmtNum
$ws.Range("E243").Value = ""

$ws.Range("F243").NumberFormat = $fmtNum
$ws.Range("F243").Value = 2385

$ws.Range("G243").NumberFormat = $fmtNum
$ws.Range("G243").Value = 34117.22

$ws.Range("H243").NumberFormat = $fmtDate
$ws.Range("H243").Value = (Get-Date -Year 2025 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)

$ws.Range("J243").HorizontalAlignment = -4152
$ws.Range("J243").Value = 'מימון אסתר          מזהה 070120555'

# Row 244
$ws.Range("A244").NumberFormat = $fmtDate
$ws.Range("A244").Value = (Get-Date -Year 2025 -Month 7 -Day 31 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B244").HorizontalAlignment = -4152
$ws.Range("B244").Value = 'האוצר-תגמולים'

$ws.Range("C244").HorizontalAlignment = -4152
$ws.Range("C244").Value = '  עבור: מימון אסתר          מזהה 070120555'

$ws.Range("D244").NumberFormat = $fmtInt
$ws.Range("D244").Value = 4061016

$ws.Range("E244").NumberFormat = $fmtNum
$ws.Range("E244").Value = ""

$ws.Range("F244").NumberFormat = $fmtNum
$ws.Range("F244").Value = 3523.15

$ws.Range("G244").NumberFormat = $fmtNum
$ws.Range("G244").Value = 39062.370000000003

$ws.Range("H244").NumberFormat = $fmtDate
$ws.Range("H244").Value = (Get-Date -Year 2025 -Month 7 -Day 31 -Hour 0 -Minute 0 -Second 0)

$ws.Range("J244").HorizontalAlignment = -4152
$ws.Range("J244").Value = 'מימון אסתר          מזהה 070120555'

# Row 245
$ws.Range("A245").NumberFormat = $fmtDate
$ws.Range("A245").Value = (Get-Date -Year 2025 -Month 8 -Day 4 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B245").HorizontalAlignment = -4152
$ws.Range("B245").Value = 'מסטרקרד'

$ws.Range("D245").NumberFormat = $fmtInt
$ws.Range("D245").Value = 961

$ws.Range("E245").NumberFormat = $fmtNum
$ws.Range("E245").Value = 4334.26

$ws.Range("F245").NumberFormat = $fmtNum
$ws.Range("F245").Value = ""

$ws.Range("G245").NumberFormat = $fmtNum
$ws.Range("G245").Value = 31990.11

$ws.Range("H245").NumberFormat = $fmtDate
$ws.Range("H245").Value = (Get-Date -Year 2025 -Month 8 -Day 4 -Hour 0 -Minute 0 -Second 0)

# Row 246
$ws.Range("A246").NumberFormat = $fmtDate
$ws.Range("A246").Value = (Get-Date -Year 2025 -Month 8 -Day 4 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B246").HorizontalAlignment = -4152
$ws.Range("B246").Value = 'העב'' לאחר-נייד'

$ws.Range("C246").HorizontalAlignment = -4152
$ws.Range("C246").Value = 'לטובת: פני דויטש עבור: שכירות אוגוסט 2025'

$ws.Range("D246").NumberFormat = $fmtInt
$ws.Range("D246").Value = 409379952

$ws.Range("E246").NumberFormat = $fmtNum
$ws.Range("E246").Value = 5300

$ws.Range("F246").NumberFormat = $fmtNum
$ws.Range("F246").Value = ""

$ws.Range("G246").NumberFormat = $fmtNum
$ws.Range("G246").Value = 36324.370000000003

$ws.Range("H246").NumberFormat = $fmtDate
$ws.Range("H246").Value = (Get-Date -Year 2025 -Month 8 -Day 4 -Hour 0 -Minute 0 -Second 0)

$ws.Range("I246").HorizontalAlignment = -4152
$ws.Range("I246").Value = 'פני דויטש'

$ws.Range("J246").HorizontalAlignment = -4152
$ws.Range("J246").Value = 'שכירות אוגוסט 2025'

# Row 247
$ws.Range("A247").NumberFormat = $fmtDate
$ws.Range("A247").Value = (Get-Date -Year 2025 -Month 8 -Day 4 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B247").HorizontalAlignment = -4152
$ws.Range("B247").Value = 'בטוח לאומי חד'

$ws.Range("C247").HorizontalAlignment = -4152
$ws.Range("C247").Value = '  עבור: מימון אסתר          מזהה 070120555'

$ws.Range("D247").NumberFormat = $fmtInt
$ws.Range("D247").Value = 63051080

$ws.Range("E247").NumberFormat = $fmtNum
$ws.Range("E247").Value = ""

$ws.Range("F247").NumberFormat = $fmtNum
$ws.Range("F247").Value = 2562

$ws.Range("G247").NumberFormat = $fmtNum
$ws.Range("G247").Value = 41624.370000000003

$ws.Range("H247").NumberFormat = $fmtDate
$ws.Range("H247").Value = (Get-Date -Year 2025 -Month 8 -Day 4 -Hour 0 -Minute 0 -Second 0)

$ws.Range("J247").HorizontalAlignment = -4152
$ws.Range("J247").Value = 'מימון אסתר          מזהה 070120555'

# Row 248
$ws.Range("A248").NumberFormat = $fmtDate
$ws.Range("A248").Value = (Get-Date -Year 2025 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B248").HorizontalAlignment = -4152
$ws.Range("B248").Value = 'ע. מסלול בסיסי'

$ws.Range("D248").NumberFormat = $fmtInt
$ws.Range("D248").Value = 817

$ws.Range("E248").NumberFormat = $fmtNum
$ws.Range("E248").Value = 10

$ws.Range("F248").NumberFormat = $fmtNum
$ws.Range("F248").Value = ""

$ws.Range("G248").NumberFormat = $fmtNum
$ws.Range("G248").Value = 35180.11

$ws.Range("H248").NumberFormat = $fmtDate
$ws.Range("H248").Value = (Get-Date -Year 2025 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)

# Row 249
$ws.Range("A249").NumberFormat = $fmtDate
$ws.Range("A249").Value = (Get-Date -Year 2025 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B249").HorizontalAlignment = -4152
$ws.Range("B249").Value = 'הפק.שיק בסלולר'

$ws.Range("D249").NumberFormat = $fmtInt
$ws.Range("D249").Value = 80000053

$ws.Range("E249").NumberFormat = $fmtNum
$ws.Range("E249").Value = ""

$ws.Range("F249").NumberFormat = $fmtNum
$ws.Range("F249").Value = 3200

$ws.Range("G249").NumberFormat = $fmtNum
$ws.Range("G249").Value = 35190.11

$ws.Range("H249").NumberFormat = $fmtDate
$ws.Range("H249").Value = (Get-Date -Year 2025 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)

# Row 250
$ws.Range("A250").NumberFormat = $fmtDate
$ws.Range("A250").Value = (Get-Date -Year 2025 -Month 8 -Day 6 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B250").HorizontalAlignment = -4152
$ws.Range("B250").Value = 'פנקסי שיקים'

$ws.Range("D250").NumberFormat = $fmtInt
$ws.Range("D250").Value = 25

$ws.Range("E250").NumberFormat = $fmtNum
$ws.Range("E250").Value = 9

$ws.Range("F250").NumberFormat = $fmtNum
$ws.Range("F250").Value = ""

$ws.Range("G250").NumberFormat = $fmtNum
$ws.Range("G250").Value = 35171.11

$ws.Range("H250").NumberFormat = $fmtDate
$ws.Range("H250").Value = (Get-Date -Year 2025 -Month 8 -Day 6 -Hour 0 -Minute 0 -Second 0)

# Row 251
$ws.Range("A251").NumberFormat = $fmtDate
$ws.Range("A251").Value = (Get-Date -Year 2025 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B251").HorizontalAlignment = -4152
$ws.Range("B251").Value = 'חברת החשמל ליש'

$ws.Range("D251").NumberFormat = $fmtInt
$ws.Range("D251").Value = 706994

$ws.Range("E251").NumberFormat = $fmtNum
$ws.Range("E251").Value = 382.75

$ws.Range("F251").NumberFormat = $fmtNum
$ws.Range("F251").Value = ""

$ws.Range("G251").NumberFormat = $fmtNum
$ws.Range("G251").Value = 31337.34

$ws.Range("H251").NumberFormat = $fmtDate
$ws.Range("H251").Value = (Get-Date -Year 2025 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Row 252
$ws.Range("A252").NumberFormat = $fmtDate
$ws.Range("A252").Value = (Get-Date -Year 2025 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B252").HorizontalAlignment = -4152
$ws.Range("B252").Value = 'שירותי בריאות'

$ws.Range("D252").NumberFormat = $fmtInt
$ws.Range("D252").Value = 30744700

$ws.Range("E252").NumberFormat = $fmtNum
$ws.Range("E252").Value = 221.02

$ws.Range("F252").NumberFormat = $fmtNum
$ws.Range("F252").Value = ""

$ws.Range("G252").NumberFormat = $fmtNum
$ws.Range("G252").Value = 31720.09

$ws.Range("H252").NumberFormat = $fmtDate
$ws.Range("H252").Value = (Get-Date -Year 2025 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Row 253
$ws.Range("A253").NumberFormat = $fmtDate
$ws.Range("A253").Value = (Get-Date -Year 2025 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B253").HorizontalAlignment = -4152
$ws.Range("B253").Value = 'העב'' לאחר-נייד'

$ws.Range("C253").HorizontalAlignment = -4152
$ws.Range("C253").Value = 'לטובת: SACHIN CHANDRAPPA עבור: salary July 2025'

$ws.Range("D253").NumberFormat = $fmtInt
$ws.Range("D253").Value = 410957145

$ws.Range("E253").NumberFormat = $fmtNum
$ws.Range("E253").Value = 3230

$ws.Range("F253").NumberFormat = $fmtNum
$ws.Range("F253").Value = ""

$ws.Range("G253").NumberFormat = $fmtNum
$ws.Range("G253").Value = 31941.11

$ws.Range("H253").NumberFormat = $fmtDate
$ws.Range("H253").Value = (Get-Date -Year 2025 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0)

$ws.Range("I253").HorizontalAlignment = -4152
$ws.Range("I253").Value = 'SACHIN CHANDRAPPA'

$ws.Range("J253").HorizontalAlignment = -4152
$ws.Range("J253").Value = 'salary July 2025'

# Row 254
$ws.Range("A254").NumberFormat = $fmtDate
$ws.Range("A254").Value = (Get-Date -Year 2025 -Month 8 -Day 15 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B254").HorizontalAlignment = -4152
$ws.Range("B254").Value = 'הו"ק פקדון'

$ws.Range("D254").NumberFormat = $fmtInt
$ws.Range("D254").Value = 2590001

$ws.Range("E254").NumberFormat = $fmtNum
$ws.Range("E254").Value = 1100

$ws.Range("F254").NumberFormat = $fmtNum
$ws.Range("F254").Value = ""

$ws.Range("G254").NumberFormat = $fmtNum
$ws.Range("G254").Value = 30237.34

$ws.Range("H254").NumberFormat = $fmtDate
$ws.Range("H254").Value = (Get-Date -Year 2025 -Month 8 -Day 15 -Hour 0 -Minute 0 -Second 0)

# Row 255
$ws.Range("A255").NumberFormat = $fmtDate
$ws.Range("A255").Value = (Get-Date -Year 2025 -Month 8 -Day 17 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B255").HorizontalAlignment = -4152
$ws.Range("B255").Value = 'סקיולייף-נטלי'

$ws.Range("D255").NumberFormat = $fmtInt
$ws.Range("D255").Value = 32349623

$ws.Range("E255").NumberFormat = $fmtNum
$ws.Range("E255").Value = 25

$ws.Range("F255").NumberFormat = $fmtNum
$ws.Range("F255").Value = ""

$ws.Range("G255").NumberFormat = $fmtNum
$ws.Range("G255").Value = 30212.34

$ws.Range("H255").NumberFormat = $fmtDate
$ws.Range("H255").Value = (Get-Date -Year 2025 -Month 8 -Day 17 -Hour 0 -Minute 0 -Second 0)

# Row 256
$ws.Range("A256").NumberFormat = $fmtDate
$ws.Range("A256").Value = (Get-Date -Year 2025 -Month 8 -Day 20 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B256").HorizontalAlignment = -4152
$ws.Range("B256").Value = 'בזק-הוראות קבע'

$ws.Range("D256").NumberFormat = $fmtInt
$ws.Range("D256").Value = 703991

$ws.Range("E256").NumberFormat = $fmtNum
$ws.Range("E256").Value = 49.16

$ws.Range("F256").NumberFormat = $fmtNum
$ws.Range("F256").Value = ""

$ws.Range("G256").NumberFormat = $fmtNum
$ws.Range("G256").Value = 30163.18

$ws.Range("H256").NumberFormat = $fmtDate
$ws.Range("H256").Value = (Get-Date -Year 2025 -Month 8 -Day 20 -Hour 0 -Minute 0 -Second 0)

# Row 257
$ws.Range("A257").NumberFormat = $fmtDate
$ws.Range("A257").Value = (Get-Date -Year 2025 -Month 8 -Day 28 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B257").HorizontalAlignment = -4152
$ws.Range("B257").Value = 'ביטוח.ל-אזו"ש'

$ws.Range("C257").HorizontalAlignment = -4152
$ws.Range("C257").Value = '  עבור: מימון אסתר          מזהה 070120555'

$ws.Range("D257").NumberFormat = $fmtInt
$ws.Range("D257").Value = 13101084

$ws.Range("E257").NumberFormat = $fmtNum
$ws.Range("E257").Value = ""

$ws.Range("F257").NumberFormat = $fmtNum
$ws.Range("F257").Value = 1422

$ws.Range("G257").NumberFormat = $fmtNum
$ws.Range("G257").Value = 33970.18

$ws.Range("H257").NumberFormat = $fmtDate
$ws.Range("H257").Value = (Get-Date -Year 2025 -Month 8 -Day 28 -Hour 0 -Minute 0 -Second 0)

$ws.Range("J257").HorizontalAlignment = -4152
$ws.Range("J257").Value = 'מימון אסתר          מזהה 070120555'

# Row 258
$ws.Range("A258").NumberFormat = $fmtDate
$ws.Range("A258").Value = (Get-Date -Year 2025 -Month 8 -Day 28 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B258").HorizontalAlignment = -4152
$ws.Range("B258").Value = 'ביטוח.ל-אזו"ש'

$ws.Range("C258").HorizontalAlignment = -4152
$ws.Range("C258").Value = '  עבור: מימון אסתר          מזהה 070120555'

$ws.Range("D258").NumberFormat = $fmtInt
$ws.Range("D258").Value = 13101084

$ws.Range("E258").NumberFormat = $fmtNum
$ws.Range("E258").Value = ""

$ws.Range("F258").NumberFormat = $fmtNum
$ws.Range("F258").Value = 2385

$ws.Range("G258").NumberFormat = $fmtNum
$ws.Range("G258").Value = 32548.18

$ws.Range("H258").NumberFormat = $fmtDate
$ws.Range("H258").Value = (Get-Date -Year 2025 -Month 8 -Day 28 -Hour 0 -Minute 0 -Second 0)

$ws.Range("J258").HorizontalAlignment = -4152
$ws.Range("J258").Value = 'מימון אסתר          מזהה 070120555'

# Row 259
$ws.Range("A259").NumberFormat = $fmtDate
$ws.Range("A259").Value = (Get-Date -Year 2025 -Month 8 -Day 31 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B259").HorizontalAlignment = -4152
$ws.Range("B259").Value = 'האוצר-תגמולים'

$ws.Range("C259").HorizontalAlignment = -4152
$ws.Range("C259").Value = '  עבור: מימון אסתר          מזהה 070120555'

$ws.Range("D259").NumberFormat = $fmtInt
$ws.Range("D259").Value = 4061016

$ws.Range("E259").NumberFormat = $fmtNum
$ws.Range("E259").Value = ""

$ws.Range("F259").NumberFormat = $fmtNum
$ws.Range("F259").Value = 3523.15

$ws.Range("G259").NumberFormat = $fmtNum
$ws.Range("G259").Value = 37493.33

$ws.Range("H259").NumberFormat = $fmtDate
$ws.Range("H259").Value = (Get-Date -Year 2025 -Month 8 -Day 31 -Hour 0 -Minute 0 -Second 0)

$ws.Range("J259").HorizontalAlignment = -4152
$ws.Range("J259").Value = 'מימון אסתר          מזהה 070120555'

# Row 260
$ws.Range("A260").NumberFormat = $fmtDate
$ws.Range("A260").Value = (Get-Date -Year 2025 -Month 9 -Day 1 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B260").HorizontalAlignment = -4152
$ws.Range("B260").Value = 'העב'' לאחר-נייד'

$ws.Range("C260").HorizontalAlignment = -4152
$ws.Range("C260").Value = 'לטובת: פני דויטש עבור: שכירות ספטמבר 2025'

$ws.Range("D260").NumberFormat = $fmtInt
$ws.Range("D260").Value = 414537948

$ws.Range("E260").NumberFormat = $fmtNum
$ws.Range("E260").Value = 5300

$ws.Range("F260").NumberFormat = $fmtNum
$ws.Range("F260").Value = ""

$ws.Range("G260").NumberFormat = $fmtNum
$ws.Range("G260").Value = 32193.33

$ws.Range("H260").NumberFormat = $fmtDate
$ws.Range("H260").Value = (Get-Date -Year 2025 -Month 9 -Day 1 -Hour 0 -Minute 0 -Second 0)

$ws.Range("I260").HorizontalAlignment = -4152
$ws.Range("I260").Value = 'פני דויטש'

$ws.Range("J260").HorizontalAlignment = -4152
$ws.Range("J260").Value = 'שכירות ספטמבר 2025'

# Row 261
$ws.Range("A261").NumberFormat = $fmtDate
$ws.Range("A261").Value = (Get-Date -Year 2025 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B261").HorizontalAlignment = -4152
$ws.Range("B261").Value = 'מסטרקרד'

$ws.Range("D261").NumberFormat = $fmtInt
$ws.Range("D261").Value = 961

$ws.Range("E261").NumberFormat = $fmtNum
$ws.Range("E261").Value = 2649.75

$ws.Range("F261").NumberFormat = $fmtNum
$ws.Range("F261").Value = ""

$ws.Range("G261").NumberFormat = $fmtNum
$ws.Range("G261").Value = 29543.58

$ws.Range("H261").NumberFormat = $fmtDate
$ws.Range("H261").Value = (Get-Date -Year 2025 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)

# Row 262
$ws.Range("A262").NumberFormat = $fmtDate
$ws.Range("A262").Value = (Get-Date -Year 2025 -Month 9 -Day 3 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B262").HorizontalAlignment = -4152
$ws.Range("B262").Value = 'ע.מפעולות-ישיר'

$ws.Range("D262").NumberFormat = $fmtInt
$ws.Range("D262").Value = 2

$ws.Range("E262").NumberFormat = $fmtNum
$ws.Range("E262").Value = 3.5

$ws.Range("F262").NumberFormat = $fmtNum
$ws.Range("F262").Value = ""

$ws.Range("G262").NumberFormat = $fmtNum
$ws.Range("G262").Value = 29530.080000000002

$ws.Range("H262").NumberFormat = $fmtDate
$ws.Range("H262").Value = (Get-Date -Year 2025 -Month 9 -Day 3 -Hour 0 -Minute 0 -Second 0)

# Row 263
$ws.Range("A263").NumberFormat = $fmtDate
$ws.Range("A263").Value = (Get-Date -Year 2025 -Month 9 -Day 3 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B263").HorizontalAlignment = -4152
$ws.Range("B263").Value = 'ע. מסלול בסיסי'

$ws.Range("D263").NumberFormat = $fmtInt
$ws.Range("D263").Value = 817

$ws.Range("E263").NumberFormat = $fmtNum
$ws.Range("E263").Value = 10

$ws.Range("F263").NumberFormat = $fmtNum
$ws.Range("F263").Value = ""

$ws.Range("G263").NumberFormat = $fmtNum
$ws.Range("G263").Value = 29533.58

$ws.Range("H263").NumberFormat = $fmtDate
$ws.Range("H263").Value = (Get-Date -Year 2025 -Month 9 -Day 3 -Hour 0 -Minute 0 -Second 0)

# Row 264
$ws.Range("A264").NumberFormat = $fmtDate
$ws.Range("A264").Value = (Get-Date -Year 2025 -Month 9 -Day 5 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B264").HorizontalAlignment = -4152
$ws.Range("B264").Value = 'הפק שיק-ע.ישיר'

$ws.Range("D264").NumberFormat = $fmtInt
$ws.Range("D264").Value = 80000054

$ws.Range("E264").NumberFormat = $fmtNum
$ws.Range("E264").Value = ""

$ws.Range("F264").NumberFormat = $fmtNum
$ws.Range("F264").Value = 3200

$ws.Range("G264").NumberFormat = $fmtNum
$ws.Range("G264").Value = 32730.080000000002

$ws.Range("H264").NumberFormat = $fmtDate
$ws.Range("H264").Value = (Get-Date -Year 2025 -Month 9 -Day 5 -Hour 0 -Minute 0 -Second 0)

# Reflect the post-edit selection/view state recorded in the workbook
$ws.Range("C213").Select()